# "add Gaussian_filter & update all the WorkTools"
# Re-saved from a different machine/Excel install (Mac path instead of the
# old Windows/OneDrive one): two toggle flags on the config sheet get turned
# on, the previously "active" H5File1 sheet hands the active/selected tab
# back to config, per-sheet selections move, the stray no-op style on the
# two header "group_path" cells (C2) is cleared, and column widths get the
# tiny re-metrics that come from Excel recalculating "best fit" on reopen.

$wb = $excel.ActiveWorkbook

$config   = $wb.Worksheets.Item("config")
$h5file1  = $wb.Worksheets.Item("H5File1")
$h5file2  = $wb.Worksheets.Item("H5File2")

# --- config sheet: flip both flags on ---------------------------------
$config.Range("B2").Value = 1
$config.Range("B3").Value = 1

# --- column width touch-ups (re-metriced "best fit" widths) -----------
$config.Range("A1").EntireColumn.ColumnWidth = 50.666666666666664

$h5file1.Range("A1").EntireColumn.ColumnWidth = 22.333333333333336
$h5file1.Range("C1:D1").EntireColumn.ColumnWidth = 32.5
$h5file1.Range("E1").EntireColumn.ColumnWidth = 87.33333333333334
$h5file1.Range("F1").EntireColumn.ColumnWidth = 10.333333333333332
$h5file1.Range("G1").EntireColumn.ColumnWidth = 11.666666666666666
$h5file1.Range("H1").EntireColumn.ColumnWidth = 23.5
$h5file1.Range("I1").EntireColumn.ColumnWidth = 19.333333333333336

$h5file2.Range("A1").EntireColumn.ColumnWidth = 22.333333333333336
$h5file2.Range("C1:D1").EntireColumn.ColumnWidth = 32.5
$h5file2.Range("E1").EntireColumn.ColumnWidth = 63.666666666666664
$h5file2.Range("F1").EntireColumn.ColumnWidth = 10.333333333333332
$h5file2.Range("G1").EntireColumn.ColumnWidth = 11.666666666666666
$h5file2.Range("H1").EntireColumn.ColumnWidth = 23.5
$h5file2.Range("I1").EntireColumn.ColumnWidth = 15.999999999999998

# --- drop the stray "applyFont only" style from the group_path header -
# (it carried no real formatting; clearing it folds the cell back onto the
# workbook default, same as the re-saved file)
$h5file1.Range("C2").ClearFormats()
$h5file2.Range("C2").ClearFormats()

# --- per-sheet selections ----------------------------------------------
$h5file2.Range("F3:G6").Select()
$h5file2.Range("F3").Activate()

$h5file1.Range("E13").Select()

$config.Range("B4").Select()

# config becomes the active/front-most tab again (was H5File1)
$config.Activate()
